$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 857, pushing the existing data
# (old rows 857..934) down to 859..936. This also grows the used range
# from A1:T934 to A1:T936.
$ws.Range("A857:A858").EntireRow.Insert()

# New row 857: new weekly price observation (Navel Late / Primera)
$ws.Cells.Item(857, 1).Value = 5
$ws.Cells.Item(857, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(857, 3).Value = "Maule"
$ws.Cells.Item(857, 4).Value = 45132
$ws.Cells.Item(857, 5).Value = 7
$ws.Cells.Item(857, 6).Value = "Fruta"
$ws.Cells.Item(857, 7).Value = 100102
$ws.Cells.Item(857, 8).Value = "Cítricos"
$ws.Cells.Item(857, 9).Value = 100102005
$ws.Cells.Item(857, 10).Value = "Naranja"
$ws.Cells.Item(857, 11).Value = "Navel Late"
$ws.Cells.Item(857, 12).Value = "Primera"
$ws.Cells.Item(857, 13).Value = 450
$ws.Cells.Item(857, 14).Value = 7000
$ws.Cells.Item(857, 15).Value = 7000
$ws.Cells.Item(857, 16).Value = 7000
$ws.Cells.Item(857, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(857, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(857, 19).Value = 467
$ws.Cells.Item(857, 20).Value = 15

# New row 858: new weekly price observation (New Hall / Primera)
$ws.Cells.Item(858, 1).Value = 5
$ws.Cells.Item(858, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(858, 3).Value = "Maule"
$ws.Cells.Item(858, 4).Value = 45132
$ws.Cells.Item(858, 5).Value = 7
$ws.Cells.Item(858, 6).Value = "Fruta"
$ws.Cells.Item(858, 7).Value = 100102
$ws.Cells.Item(858, 8).Value = "Cítricos"
$ws.Cells.Item(858, 9).Value = 100102005
$ws.Cells.Item(858, 10).Value = "Naranja"
$ws.Cells.Item(858, 11).Value = "New Hall"
$ws.Cells.Item(858, 12).Value = "Primera"
$ws.Cells.Item(858, 13).Value = 300
$ws.Cells.Item(858, 14).Value = 7000
$ws.Cells.Item(858, 15).Value = 7000
$ws.Cells.Item(858, 16).Value = 7000
$ws.Cells.Item(858, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(858, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(858, 19).Value = 467
$ws.Cells.Item(858, 20).Value = 15
